$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schedule")

$ws.Range("D8").Value = '* Describe the difference between training and testing data sets
* Describe the differnece between a parametric and non-parametric model
* Identify  and describe situations where classification, regression, and clustering models are appropriate. 
* Explain the concept of overfitting, and bias-variance tradeoff. '
$ws.Range("E8").Value = 'Summarize chosen DDS Case study on Tue  
Read ISLR Ch1 before Th. 
Prepare your part of HW5 before Th.'
$ws.Range("F8").Value = 'ISLR Ch 1, 2  
ISLR MOOC [Video lectures](http://auapps.american.edu/alberto/www/analytics/ISLRLectures.html)  
[Materials from DSO/IOM 530 at USC](https://www.alsharif.info/iom530)  '
$ws.Range("G8").Value = 'Review POLS Learn R HW 3 (Due 3/3 )  
Blog post & response (Due 2/29 )  
[hw05: Statistical Learning](hw/hw05-statistical-learning.html) (Due x/xx )   '
$ws.Range("H8").Value = 'Share DDS Case study
Return to bridge (time pending)'
$ws.Range("I8").Value = 'Jigsaw lecture on ISLR Ch 2 by answering assigned questions from ISLR'
$ws.Range("C9").Value = 'Prediction using Regression Models'
$ws.Range("D9").Value = '* Translate a business question into a mathematical model. 
* Identify cases where a non-parametric approach would be more appropriate than a parametric approach. '
$ws.Range("E9").Value = 'Familarize yourself with the advertising data from Ch 2. '
$ws.Range("F9").Value = 'ISLR Ch 3. '
$ws.Range("H9").Value = 'Finish Part I of HW 5'
$ws.Range("I9").Value = 'Part III hw 5'
